$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 550, shifting existing rows 550-570 down to 551-571.
$ws.Rows.Item(550).Insert()

# Populate the newly inserted row 550 with the new weekly record.
$ws.Cells.Item(550, 1).Value = 3
$ws.Cells.Item(550, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(550, 3).Value = "Coquimbo"
$ws.Cells.Item(550, 4).Value = 45075
$ws.Cells.Item(550, 5).Value = 5
$ws.Cells.Item(550, 6).Value = 100112043
$ws.Cells.Item(550, 7).Value = "Pepino ensalada"
$ws.Cells.Item(550, 8).Value = "Sin especificar"
$ws.Cells.Item(550, 9).Value = "Primera"
$ws.Cells.Item(550, 10).Value = 105
$ws.Cells.Item(550, 11).Value = 9500
$ws.Cells.Item(550, 12).Value = 10000
$ws.Cells.Item(550, 13).Value = 9738
$ws.Cells.Item(550, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(550, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(550, 16).Value = 162
$ws.Cells.Item(550, 17).Value = 60
$ws.Cells.Item(550, 18).Value = "Hortaliza"
